$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value for column F (dSF) updates
$updates = @{
    "F6"  = -7
    "F10" = 0
    "F14" = -4
    "F16" = -1
    "F17" = -1
    "F19" = 0
    "F23" = 1
    "F31" = 0
    "F34" = -2
    "F35" = 6
    "F36" = -5
    "F38" = -2
    "F40" = 3
    "F41" = -7
    "F42" = -1
    "F44" = 0
    "F45" = 5
    "F47" = -3
    "F49" = -10
    "F54" = 2
    "F60" = -2
    "F66" = 6
    "F68" = 0
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
